$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores 16 "year" groups of 4 quarter rows (A/B/C/D).
# In each group the rows labelled "B" and "C" need to swap places
# (all of columns A-E), matching the target revision.

# Swap rows 3 and 4 (2000年B <-> 2000年C)
$ws.Range("A3").Value = "2000年C"
$ws.Range("A4").Value = "2000年B"
$ws.Range("D3").Value = 71.5
$ws.Range("D4").Value = 42
$ws.Range("E3").Value = 519.4
$ws.Range("E4").Value = 315.4

# Swap rows 7 and 8 (2001年B <-> 2001年C)
$ws.Range("A7").Value = "2001年C"
$ws.Range("A8").Value = "2001年B"
$ws.Range("B7").Value = 100
$ws.Range("B8").Value = 99.40000000000001
$ws.Range("C7").Value = 2
$ws.Range("C8").Value = ""
$ws.Range("D7").Value = 2.1
$ws.Range("D8").Value = 3.4
$ws.Range("E7").Value = 593.1
$ws.Range("E8").Value = 345.8

# Swap rows 11 and 12 (2002年B <-> 2002年C)
$ws.Range("A11").Value = "2002年C"
$ws.Range("A12").Value = "2002年B"
$ws.Range("B11").Value = 99.40000000000001
$ws.Range("B12").Value = 98.8
$ws.Range("D11").Value = 33.7
$ws.Range("D12").Value = 35.3
$ws.Range("E11").Value = 979.5
$ws.Range("E12").Value = 551.3

# Swap rows 15 and 16 (2003年B <-> 2003年C)
$ws.Range("A15").Value = "2003年C"
$ws.Range("A16").Value = "2003年B"
$ws.Range("B15").Value = 98.09999999999999
$ws.Range("B16").Value = 99.59999999999999
$ws.Range("C15").Value = -1.3
$ws.Range("C16").Value = 0.8
$ws.Range("D15").Value = 95.90000000000001
$ws.Range("D16").Value = 38.8
$ws.Range("E15").Value = 2060.6
$ws.Range("E16").Value = 1332.6

# Swap rows 19 and 20 (2004年B <-> 2004年C)
$ws.Range("A19").Value = "2004年C"
$ws.Range("A20").Value = "2004年B"
$ws.Range("B19").Value = 95.90000000000001
$ws.Range("B20").Value = 98.90000000000001
$ws.Range("C19").Value = -2.2
$ws.Range("C20").Value = -0.7
$ws.Range("D19").Value = 343.4
$ws.Range("D20").Value = 101.8
$ws.Range("E19").Value = 2894
$ws.Range("E20").Value = 1759.7

# Swap rows 23 and 24 (2005年B <-> 2005年C)
$ws.Range("A23").Value = "2005年C"
$ws.Range("A24").Value = "2005年B"
$ws.Range("B23").Value = 99.8
$ws.Range("B24").Value = 100.2
$ws.Range("C23").Value = 3.9
$ws.Range("C24").Value = 1.3
$ws.Range("D23").Value = 7.7
$ws.Range("D24").Value = -8.6
$ws.Range("E23").Value = 5540.5
$ws.Range("E24").Value = 3446.9

# Swap rows 27 and 28 (2006年B <-> 2006年C)
$ws.Range("A27").Value = "2006年C"
$ws.Range("A28").Value = "2006年B"
$ws.Range("B27").Value = 100
$ws.Range("B28").Value = 99.40000000000001
$ws.Range("C27").Value = 0.2
$ws.Range("C28").Value = -0.8
$ws.Range("D27").Value = -3.5
$ws.Range("D28").Value = 29.3
$ws.Range("E27").Value = 7053.4
$ws.Range("E28").Value = 4463.8

# Swap rows 31 and 32 (2007年B <-> 2007年C)
$ws.Range("A31").Value = "2007年C"
$ws.Range("A32").Value = "2007年B"
$ws.Range("B31").Value = 98.90000000000001
$ws.Range("B32").Value = 99.8
$ws.Range("C31").Value = -1.1
$ws.Range("C32").Value = 0.4
$ws.Range("D31").Value = 129.9
$ws.Range("D32").Value = 15.3
$ws.Range("E31").Value = 9746.1
$ws.Range("E32").Value = 5883.4

# Swap rows 35 and 36 (2008年B <-> 2008年C)
$ws.Range("A35").Value = "2008年C"
$ws.Range("A36").Value = "2008年B"
$ws.Range("B35").Value = 97.59999999999999
$ws.Range("B36").Value = 96.40000000000001
$ws.Range("C35").Value = 0.2
$ws.Range("C36").Value = -4.3
$ws.Range("D35").Value = 144.3
$ws.Range("D36").Value = 137.7
$ws.Range("E35").Value = 10230.3
$ws.Range("E36").Value = 6406.6

# Swap rows 39 and 40 (2009年B <-> 2009年C)
$ws.Range("A39").Value = "2009年C"
$ws.Range("A40").Value = "2009年B"
$ws.Range("C39").Value = 2.8
$ws.Range("C40").Value = 3.7
$ws.Range("D39").Value = -3.1
$ws.Range("D40").Value = 0.8
$ws.Range("E39").Value = 12253.8
$ws.Range("E40").Value = 7230.3

# Swap rows 43 and 44 (2010年B <-> 2010年C)
$ws.Range("A43").Value = "2010年C"
$ws.Range("A44").Value = "2010年B"
$ws.Range("B43").Value = 99.59999999999999
$ws.Range("B44").Value = 100.1
$ws.Range("C43").Value = -0.4
$ws.Range("C44").Value = 0.2
$ws.Range("D43").Value = 18
$ws.Range("D44").Value = 8.199999999999999
$ws.Range("E43").Value = 17296.2
$ws.Range("E44").Value = 10846

# Swap rows 47 and 48 (2011年B <-> 2011年C)
$ws.Range("A47").Value = "2011年C"
$ws.Range("A48").Value = "2011年B"
$ws.Range("B47").Value = 99.2
$ws.Range("B48").Value = 99.09999999999999
$ws.Range("C47").Value = -0.4
$ws.Range("C48").Value = -0.2
$ws.Range("D47").Value = 39.1
$ws.Range("D48").Value = 26.4
$ws.Range("E47").Value = 23339
$ws.Range("E48").Value = 14350.2

# Swap rows 51 and 52 (2012年B <-> 2012年C)
$ws.Range("A51").Value = "2012年C"
$ws.Range("A52").Value = "2012年B"
$ws.Range("C51").Value = 0.4
$ws.Range("C52").Value = 0.3
$ws.Range("D51").Value = 17.7
$ws.Range("D52").Value = 13.6
$ws.Range("E51").Value = 23875.3
$ws.Range("E52").Value = 15663.9

# Swap rows 55 and 56 (2013年B <-> 2013年C)
$ws.Range("A55").Value = "2013年C"
$ws.Range("A56").Value = "2013年B"
$ws.Range("B55").Value = 99.8
$ws.Range("B56").Value = 100.3
$ws.Range("C55").Value = 0.3
$ws.Range("C56").Value = -0.2
$ws.Range("D55").Value = 5.2
$ws.Range("D56").Value = 5.8
$ws.Range("E55").Value = 23744.3
$ws.Range("E56").Value = 15140

# Swap rows 59 and 60 (2014年B <-> 2014年C)
$ws.Range("A59").Value = "2014年C"
$ws.Range("A60").Value = "2014年B"
$ws.Range("B59").Value = 99.90000000000001
$ws.Range("B60").Value = 99.8
$ws.Range("C59").Value = 0.6
$ws.Range("C60").Value = -0.5
$ws.Range("D59").Value = 10.3
$ws.Range("D60").Value = 9.300000000000001
$ws.Range("E59").Value = 24378.1
$ws.Range("E60").Value = 16044.7

# Swap rows 63 and 64 (2015年B <-> 2015年C)
$ws.Range("A63").Value = "2015年C"
$ws.Range("A64").Value = "2015年B"
$ws.Range("B63").Value = 99.90000000000001
$ws.Range("B64").Value = 99.2
$ws.Range("C63").Value = 1.1
$ws.Range("C64").Value = 0.6
$ws.Range("D63").Value = -8.6
$ws.Range("D64").Value = -6.8
$ws.Range("E63").Value = 21806.7
$ws.Range("E64").Value = 14337.4

# Remove columns F and G entirely (also updates the sheet dimension)
$ws.Columns("F:G").Delete()
